# "Graba Lat Long en Entregados"
# Adds two new pending-request rows to the "Hoja1" tracking sheet:
#   Row 9  -> "Que salga OrderId al lado de N° de Remito"
#   Row 10 -> "Que cuando graba el envío, si es Entregado, grabe las
#              coordenadas en nuevos campos Latitud2 y Longitud2"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Write row 10's text first so it is registered in the shared-string
# table ahead of row 9's text (matches the order the strings were
# originally authored in: "Que cuando graba..." got the lower index).
$ws.Range("B10").Value = "Que cuando graba el envío, si es Entregado, grabe las coordenadas en nuevos campos Latitud2 y Longitud2"
$ws.Range("B9").Value = "Que salga OrderId al lado de N° de Remito"

# Row 9: fecha solicitada 5/4/2022, subida a la tienda 7/4/2022
$ws.Range("A9").Value = (Get-Date -Year 2022 -Month 4 -Day 5).Date
$ws.Range("C9").Value = (Get-Date -Year 2022 -Month 4 -Day 7).Date

# Row 10: fecha solicitada 5/4/2022, subida a la tienda 7/4/2022
$ws.Range("A10").Value = (Get-Date -Year 2022 -Month 4 -Day 5).Date
$ws.Range("C10").Value = (Get-Date -Year 2022 -Month 4 -Day 7).Date

# Row 10 grows taller so the long wrapped text fits.
$ws.Rows.Item(10).RowHeight = 20.4

# Leave the selection where the user ended up after typing the new row.
$ws.Range("B10").Select()
